# Updated Schedule Source to support also .data and .csv files.
# Files format is updated to support label creation on items from file:
# the sheet now stores a Time/Name/Q/type table instead of the old
# item1..item4 sample rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row
$ws.Range("A1").Value = "Time"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Q"
$ws.Range("D1").Value = "type"

# New data rows (2-5) - four "Plancha" entries labelled A-D
$data = @(
    @(0, "Plancha", 1, "A"),
    @(0, "Plancha", 1, "B"),
    @(0, "Plancha", 1, "C"),
    @(0, "Plancha", 10, "D")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# The lone formatted/empty marker cell moves from H14 down to H15
$ws.Rows("14").Delete()
$ws.Range("H15").Font.Underline = $true

# Selection moved as part of the edit
$ws.Range("F13").Select()

# Page setup now explicit (A4, portrait) instead of inherited defaults
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
